$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values (e.g. "602.65")
# are not auto-converted into actual numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2: 'Bitcoin'
$ws.Range("D2").Value = '63.366.39'
$ws.Range("E2").Value = '  +1.53%  '

# Row 3: 'Ethereum'
$ws.Range("D3").Value = '3.169.42'
$ws.Range("E3").Value = '  -0.64%  '

# Row 4: 'TetherUSD'
$ws.Range("E4").Value = '  -0.13%  '

# Row 5: 'BNB'
$ws.Range("D5").Value = '602.65'
$ws.Range("E5").Value = '  +1.90%  '

# Row 6: 'Solana'
$ws.Range("D6").Value = '135.68'
$ws.Range("E6").Value = '  -0.16%  '

# Row 7: 'USDC'
$ws.Range("E7").Value = '  -0.09%  '

# Row 8: 'LidoStakedEther'
$ws.Range("D8").Value = '3.168.55'
$ws.Range("E8").Value = '  -0.51%  '

# Row 9: 'XRP'
$ws.Range("D9").Value = '0.513'
$ws.Range("E9").Value = '  +1.89%  '

# Row 10: 'Dogecoin'
$ws.Range("D10").Value = '0.143'
$ws.Range("E10").Value = '  -0.35%  '

# Row 11: 'Toncoin'
$ws.Range("D11").Value = '5.34'
$ws.Range("E11").Value = '  +1.29%  '

# Row 12: 'Cardano'
$ws.Range("D12").Value = '0.454'
$ws.Range("E12").Value = '  +0.53%  '

# Row 13: 'ShibaInu'
$ws.Range("D13").Value = '0.0000240'
$ws.Range("E13").Value = '  +1.54%  '

# Row 14: 'Avalanche'
$ws.Range("D14").Value = '34.83'
$ws.Range("E14").Value = '  +4.37%  '

# Row 15: 'WrappedliquidstakedEther2.0'
$ws.Range("D15").Value = '3.690.40'
$ws.Range("E15").Value = '  -0.48%  '

# Row 16: 'TRON'
$ws.Range("E16").Value = '  +0.11%  '

# Row 17: 'WrappedEther'
$ws.Range("D17").Value = '3.170.39'
$ws.Range("E17").Value = '  -0.80%  '

# Row 18: 'WrappedBTC'
$ws.Range("D18").Value = '63.332.41'
$ws.Range("E18").Value = '  +1.40%  '

# Row 19: 'Polkadot'
$ws.Range("E19").Value = '  -0.42%  '

# Row 20: 'BitcoinCash'
$ws.Range("D20").Value = '461.70'
$ws.Range("E20").Value = '  +0.19%  '

# Row 21: 'Chainlink'
$ws.Range("E21").Value = '  +0.53%  '

# Row 22: 'Polygon'
$ws.Range("D22").Value = '0.697'
$ws.Range("E22").Value = '  -1.48%  '

# Row 23: 'Uniswap'
$ws.Range("D23").Value = '7.66'
$ws.Range("E23").Value = '  +0.43%  '

# Row 24: 'InternetComputer(DFINITY)'
$ws.Range("D24").Value = '13.23'
$ws.Range("E24").Value = '  -1.10%  '

# Row 25: 'Litecoin'
$ws.Range("D25").Value = '83.12'
$ws.Range("E25").Value = '  +0.16%  '

# Row 26: 'Dai'
$ws.Range("E26").Value = '  +0.07%  '

# Row 27: 'PancakeSwap'
$ws.Range("D27").Value = '2.70'
$ws.Range("E27").Value = '  +0.32%  '

# Row 28: 'FirstDigitalUSD'
$ws.Range("E28").Value = '  -0.05%  '

# Row 29: 'ImmutableX'
$ws.Range("E29").Value = '  +1.11%  '

# Row 30: 'NEARProtocol' -> 'RenderToken'
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '7.72'
$ws.Range("E30").Value = '  -1.99%  '

# Row 31: 'RenderToken' -> 'NEARProtocol'
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = '6.78'
$ws.Range("E31").Value = '  -1.83%  '

# Row 32: 'EthereumClassic'
$ws.Range("D32").Value = '27.18'
$ws.Range("E32").Value = '  -0.38%  '

# Row 33: 'Hedera'
$ws.Range("D33").Value = '0.101'
$ws.Range("E33").Value = '  -1.93%  '

# Row 34: 'Stacks'
$ws.Range("E34").Value = '  +0.62%  '

# Row 35: 'Mantle'
$ws.Range("E35").Value = '  -2.34%  '

# Row 36: 'Filecoin'
$ws.Range("D36").Value = '5.89'
$ws.Range("E36").Value = '  +1.08%  '

# Row 37: 'PEPE' -> 'OKB'
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Value = '51.44'
$ws.Range("E37").Value = '  +0.67%  '

# Row 38: 'OKB' -> 'PEPE'
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = '0.0₃0732'
$ws.Range("E38").Value = '  +5.62%  '

# Row 39: 'VeChain'
$ws.Range("E39").Value = '  +0.69%  '

# Row 40: 'Cosmos'
$ws.Range("D40").Value = '8.14'
$ws.Range("E40").Value = '  +1.36%  '

# Row 41: 'Kaspa'
$ws.Range("E41").Value = '  -0.99%  '

# Row 42: 'Bittensor' -> 'dogwifhat'
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = '2.64'
$ws.Range("E42").Value = '  -1.37%  '

# Row 43: 'dogwifhat' -> 'Bittensor'
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = '392.81'
$ws.Range("E43").Value = '  -4.60%  '

# Row 44: 'Maker'
$ws.Range("D44").Value = '2.793.15'
$ws.Range("E44").Value = '  -6.13%  '

# Row 45: 'TheGraph'
$ws.Range("D45").Value = '0.250'
$ws.Range("E45").Value = '  +0.03%  '

# Row 46: 'USDe' -> 'Arweave'
$ws.Range("B46").Value = 'Arweave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D46").Value = '35.86'
$ws.Range("E46").Value = '  -0.09%  '

# Row 47: 'Monero' -> 'USDe'
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").Value = '0.999'
$ws.Range("E47").Value = '  +0.01%  '

# Row 48: 'Arweave' -> 'Monero'
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = '126.81'
$ws.Range("E48").Value = '  +2.25%  '

# Row 49: 'Fetch.AI'
$ws.Range("E49").Value = '  -1.88%  '

# Row 50: 'InjectiveProtocol'
$ws.Range("D50").Value = '25.11'
$ws.Range("E50").Value = '  -2.28%  '

# Row 51: 'Stellar'
$ws.Range("D51").Value = '0.111'
$ws.Range("E51").Value = '  +0.18%  '
